$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update time_slot_availability values to ensure no two consecutive timeslots
# in the same day occur for these invigilators (constraint added), and fix
# the lead flag for Abigail Smith.
$ws.Range("C2").Value = "1,2,3,4,6,7,9,10,12,13,15,16,18"
$ws.Range("D2").Value = 1
$ws.Range("C3").Value = "1,2,4,5,7,8,10,11,13,14"
$ws.Range("C6").Value = "2,4,7,8,10,11,13,15,16,17,18"
$ws.Range("C5").Value = "1,2,3,6,7,8,9,10,12,15,16,17,18"

# Update the active selection on the sheet to reflect the last-edited cell.
$ws.Range("C5").Select()

# Update the workbook window position to reflect the saved window state.
$excel.ActiveWindow.WindowState = -4143
$excel.Left = -1260
$excel.Top = 1040
